# "added upload only endpoint" - update the ACI-Bootstrap-Tool sample row
# with demo/POV values instead of the old test values.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ACI-Bootstrap-Tool")

# Form Date Field: 05/03/17 -> 05/15/17 (serial 42858 -> 42870)
$ws.Range("A2").Value = "05/15/17"

# customer_name
$ws.Range("C2").Value = "demo-customer"

# router_speed (was empty) / l2-out-switch-speed
$ws.Range("M2").Value = "1G-auto"
$ws.Range("N2").Value = "1G-auto"

# switch_oob_pool_end
$ws.Range("P2").Value = "10.94.238.160"

# vcenter_ip / vcenter_datacenter / vcenter_password
$ws.Range("T2").Value = "10.94.238.58"
$ws.Range("U2").Value = "demo-dc"
$ws.Range("W2").Value = "C!sco123"

# tenant_name_1
$ws.Range("X2").Value = "POV-3tier-L3_1"

# network definitions for tenant 1
$ws.Range("Z2").Value = "10.94.238.1/24"
$ws.Range("AA2").Value = "private"
$ws.Range("AB2").Value = "172.16.32.1/24"
$ws.Range("AD2").Value = "172.16.33.1/24"
$ws.Range("AF2").Value = "172.16.34.1/24"

# ospf_id_1
$ws.Range("AH2").Value = 10

# eigrp_pid_1 cleared (protocol now ospf for tenant 1)
$ws.Range("AJ2").ClearContents()

# tenant_name_2 cleared
$ws.Range("AN2").ClearContents()

# routed_out_protocol_2: eigrp -> static
$ws.Range("AW2").Value = "static"

# eigrp_pid_2 cleared
$ws.Range("AZ2").ClearContents()

# leaf201_routed_addr_2 / leaf202_routed_addr_2 cleared
$ws.Range("BB2").ClearContents()
$ws.Range("BC2").ClearContents()
